# Update the date heading.
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2026-02-10 Tuesday"

# Update the division-problem answers in the single table, cell by cell
# (direct Range.Text assignment keeps this scoped to just the one cell,
# unlike Find/Replace which here operates document-wide).
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "30÷7=4, 2"
$t.Cell(1,2).Range.Text  = "83÷3=27, 2"
$t.Cell(1,3).Range.Text  = "61÷4=15, 1"
$t.Cell(1,4).Range.Text  = "54÷4=13, 2"
$t.Cell(1,5).Range.Text  = "89÷4=22, 1"

$t.Cell(5,1).Range.Text  = "11÷9=1, 2"
$t.Cell(5,2).Range.Text  = "83÷7=11, 6"
$t.Cell(5,3).Range.Text  = "84÷7=12, 0"
$t.Cell(5,4).Range.Text  = "64÷4=16, 0"
$t.Cell(5,5).Range.Text  = "65÷7=9, 2"

$t.Cell(9,1).Range.Text  = "34÷9=3, 7"
$t.Cell(9,2).Range.Text  = "57÷9=6, 3"
$t.Cell(9,3).Range.Text  = "18÷8=2, 2"
$t.Cell(9,4).Range.Text  = "27÷8=3, 3"
$t.Cell(9,5).Range.Text  = "89÷2=44, 1"

$t.Cell(13,1).Range.Text = "17÷8=2, 1"
$t.Cell(13,2).Range.Text = "90÷8=11, 2"
$t.Cell(13,3).Range.Text = "60÷2=30, 0"
$t.Cell(13,4).Range.Text = "55÷7=7, 6"
$t.Cell(13,5).Range.Text = "26÷8=3, 2"

$t.Cell(17,1).Range.Text = "12÷2=6, 0"
$t.Cell(17,2).Range.Text = "11÷5=2, 1"
$t.Cell(17,3).Range.Text = "17÷5=3, 2"
$t.Cell(17,4).Range.Text = "55÷5=11, 0"
$t.Cell(17,5).Range.Text = "42÷4=10, 2"
